# Auto-generated Excel COM-interop script applying the crypto price/volume diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.408.96"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "2.884.45"
$ws.Range("E3").Value = "  -3.61%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'518.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.65%  "
$ws.Range("D6").Value = "'141.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.09%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.81%  "
$ws.Range("D9").Value = "2.892.42"
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'6.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.95%  "
$ws.Range("D12").Value = "'0.354"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("D13").Value = "3.393.27"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "60.554.72"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "'22.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.62%  "
$ws.Range("D17").Value = "2.895.55"
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").Value = "'0.0000139"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.90%  "
$ws.Range("D19").Value = "'4.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.49%  "
$ws.Range("D20").Value = "'11.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.64%  "
$ws.Range("D21").Value = "'349.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.97%  "
$ws.Range("D22").Value = "'6.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.16%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'64.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").Value = "'0.444"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.54%  "
$ws.Range("D27").Value = "'0.176"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.31%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "'7.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.02%  "
$ws.Range("D30").Value = "0.0₃0849"
$ws.Range("E30").Value = "  -9.14%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'1.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("D33").Value = "'19.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("D34").Value = "'153.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("D35").Value = "'4.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.18%  "
$ws.Range("D36").Value = "'5.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.78%  "
$ws.Range("D37").Value = "'0.975"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.82%  "
$ws.Range("D38").Value = "'1.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.84%  "
$ws.Range("D39").Value = "'37.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "'1.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.66%  "
$ws.Range("D41").Value = "'0.647"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("D42").Value = "2.264.65"
$ws.Range("E42").Value = "  -6.35%  "
$ws.Range("D43").Value = "'3.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.62%  "
$ws.Range("D44").Value = "'0.0574"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.09%  "
$ws.Range("D45").Value = "'20.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.98%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'4.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.54%  "
$ws.Range("D48").Value = "'0.0233"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("D49").Value = "'10.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "'0.0904"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.30%  "
$ws.Range("D51").Value = "'18.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.65%  "
